# bitaxeMax BOM.xlsx edit script
# Commit: "added manufacturing files. clarified MCP1824 part names"
#
# Semantic changes applied:
#  1. BOM sheet: JP4 (row 30) marked DNP ("T" in column F), matching the
#     neighboring jumper/solder-jumper rows that are also not populated.
#  2. DK Order sheet: the capacitor reference-designator list in row 3 is
#     rewritten using the shorter "ranged" notation
#     ("C2-6, C8, C10-12, C14-18, ...") instead of the fully spelled out
#     list of individual refs.
#  3. DK Order sheet: JP4 has no Digikey order part (it's just a solder
#     jumper), so its row is removed from the purchasing list entirely;
#     everything below shifts up by one row.
#  4. The "DK Order" sheet becomes the active/selected sheet, with the
#     data block selected.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("BOM")
$ws2 = $wb.Worksheets.Item("DK Order")

# --- 1. BOM sheet: mark JP4 (row 30) as DNP ---------------------------------
$ws1.Range("F30").Value = "T"

# --- 2. DK Order sheet: consolidate the capacitor ref-des list -------------
$ws2.Range("B3").Value = "C2-6, C8, C10-12, C14-18, C24, C25, C27, C35, C44, C49, C50"

# --- 3. DK Order sheet: remove the JP4 row (no purchasable part) -----------
$ws2.Rows(20).Delete()

# --- 4. Make "DK Order" the active sheet with its data selected ------------
$ws2.Activate() | Out-Null
$ws2.Range("A2:C48").Select() | Out-Null
